# Update 'Experience' column (J) values on Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("J3").Value = 2
$ws.Range("J6").Value = 4
$ws.Range("J7").Value = 3
$ws.Range("J8").Value = 2
$ws.Range("J9").Value = 3
$ws.Range("J10").Value = 1
